$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently has 2 worker/debt rows (16 & 17). The new account
# statement has 13 rows covering 6 workers x up to 3 periods each.
# Insert 11 new rows after row 16 (before the old row 17), which pushes the
# old row 17 down to row 28 and the signature rows (22/23) down to (33/34).
# ---------------------------------------------------------------------------
$ws.Rows("17:27").Insert()

# Copy the formatting of the template data row (16) onto the newly inserted
# rows so they match the look of the rest of the table.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J27").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Update the summary figures at the top of the statement.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 719342      # VALOR MORA total
$ws.Range("C13").Value = 6           # Cant. Trabajadores
$ws.Range("F13").Value = 3           # Cant. Periodos

# ---------------------------------------------------------------------------
# Fill in the worker/debt detail rows (16 through 28).
# Columns: B = Tipo Doc, C = N Doc Trabajador, D = Nombre Trabajador,
#          E = Periodo Mora, F = Valor Mora, G = Salario Basico
# ---------------------------------------------------------------------------
$data = @(
    @("CC", "1047443992", "JOSE ARLY GOMEZ AGRESOT", "2507", 56940, 1160000),
    @("CC", "1047443992", "JOSE ARLY GOMEZ AGRESOT", "2506", 56940, 1160000),
    @("CC", "1047443992", "JOSE ARLY GOMEZ AGRESOT", "2504", 55042, 1160000),
    @("CC", "9237021", "DONICEL POLO SUAREZ", "2507", 56940, 1423500),
    @("CC", "9237021", "DONICEL POLO SUAREZ", "2506", 56940, 1423500),
    @("CC", "11165502", "JULIO CESAR MORELO MACHADO", "2507", 56940, 1423500),
    @("CC", "11165502", "JULIO CESAR MORELO MACHADO", "2506", 56940, 1423500),
    @("CC", "1151202337", "LEYDER ALBERTO CARMONA DE ARCO", "2507", 56940, 1423500),
    @("CC", "1151202337", "LEYDER ALBERTO CARMONA DE ARCO", "2506", 56940, 1423500),
    @("CC", "1002245176", "GLADYS SARAY HERNANDEZ MORA", "2507", 56940, 1423500),
    @("CC", "1002245176", "GLADYS SARAY HERNANDEZ MORA", "2506", 56940, 1423500),
    @("CC", "1093925873", "MAIKET JESUS AGUILAR PEÃ?A", "2507", 56940, 1423500),
    @("CC", "1093925873", "MAIKET JESUS AGUILAR PEÃ?A", "2506", 37960, 1423500)
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
}
